$wb = $excel.ActiveWorkbook
$ws = $wb.ActiveSheet

# The sheet currently has data in rows 1-28 (A1:G28). Append a new row (29)
# with a new vessel entry, matching the layout/format of the row above it.
$lastRow = 28
$newRow = $lastRow + 1

# Copy the row above so the new row inherits the same styling
# (centered text, borders, date format) as the rest of the table.
$ws.Range("B" + $lastRow + ":G" + $lastRow).Copy()
$ws.Range("B" + $newRow + ":G" + $newRow).PasteSpecial(-4122)

$ws.Cells.Item($newRow, 2).Value = "EP TASA 17 ARBOLADURA"
$ws.Cells.Item($newRow, 3).Value = "Embarcación Pesquera"
$ws.Cells.Item($newRow, 4).Value = "GP/75"
$ws.Cells.Item($newRow, 5).Value = "GP/75-224"
$ws.Cells.Item($newRow, 6).Value = (Get-Date -Year 2024 -Month 12 -Day 12).Date
$ws.Cells.Item($newRow, 7).Value = (Get-Date -Year 2024 -Month 12 -Day 20).Date

$excel.ActiveWindow.ScrollRow = 3
$excel.ActiveWindow.ScrollColumn = 2
$ws.Range("G30").Select()
